# Send Email to Business
# Adds EmailBody / EmailTo / EmailSubject asset rows to the "Assets" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Assets")

$ws.Range("B8").Value = "ADYF.EmailBody"
$ws.Range("A8").Value = "EmailBody"
$ws.Range("A9").Value = "EmailTo"
$ws.Range("A10").Value = "EmailSubject"
$ws.Range("B9").Value = "ADYF.EmailTo"
$ws.Range("B10").Value = "ADYF.EmailSubject"
$ws.Range("C8").Value = "Email Body"
$ws.Range("C9").Value = "Email Recepient"
$ws.Range("C10").Value = "Email Subject"

$ws.Range("A8:C10").Style = $ws.Range("A7:C7").Style

$ws.Activate()
$ws.Range("C10").Select()
